$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($cell, [string]$val)
    # Prefix with an apostrophe so Excel always stores the value as literal
    # text (matching the source workbook's inlineStr / text cells), even
    # when the value looks like a number. Resetting the style afterwards
    # clears the "quote prefix" cell-format flag that the apostrophe
    # entry-mode would otherwise leave behind, so formatting matches the
    # original (default/general) style.
    $cell.Value = "'" + $val
    $cell.Style = "Normal"
}

# Column D (Price) and Column E (Volume(1h)) updates per row.
# Each entry: row, newD (or $null if unchanged), newE (or $null if unchanged)
$updates = @(
    @(2,  "64.672.21",  "  +1.38%  "),
    @(3,  "3.445.63",   "  +1.09%  "),
    @(4,  $null,        "  +0.00%  "),
    @(5,  "573.94",     "  +0.96%  "),
    @(6,  "160.03",     "  +2.23%  "),
    @(7,  $null,        "  +15.01%  "),
    @(8,  $null,        "  +0.08%  "),
    @(9,  "3.446.47",   "  +1.09%  "),
    @(10, "7.22",       "  -1.65%  "),
    @(11, $null,        "  +1.85%  "),
    @(12, "0.447",      "  +3.34%  "),
    @(13, "4.044.86",   "  +1.38%  "),
    @(14, $null,        "  +0.65%  "),
    @(15, "0.0000191",  "  +0.27%  "),
    @(16, "28.06",      "  +3.13%  "),
    @(17, "64.716.88",  "  +1.51%  "),
    @(18, "3.422.90",   "  +0.41%  "),
    @(19, "6.43",       "  +2.83%  "),
    @(20, "14.27",      "  +1.84%  "),
    @(21, "379.13",     "  -0.51%  "),
    @(22, "8.04",       "  -0.07%  "),
    @(23, "0.549",      "  +3.33%  "),
    @(24, "72.70",      "  +1.15%  "),
    @(25, "0.998",      "  -0.22%  "),
    @(26, "0.0000118",  "  -1.05%  "),
    @(27, "9.92",       "  +5.50%  "),
    @(28, "0.177",      "  -0.62%  "),
    @(29, $null,        "  -0.15%  "),
    @(30, "1.50",       "  +10.23%  "),
    @(31, "6.08",       "  -1.14%  "),
    @(32, "2.03",       "  +1.70%  "),
    @(33, "23.48",      "  +0.99%  "),
    @(34, "7.17",       "  +5.44%  "),
    @(35, "1.60",       "  +10.00%  "),
    @(36, "161.41",     "  +0.84%  "),
    @(37, "1.91",       "  +4.25%  "),
    @(38, "0.0770",     "  +1.92%  "),
    @(39, "2.942.71",   "  -1.26%  "),
    @(40, "26.41",      "  -1.60%  "),
    @(41, "4.56",       "  +6.16%  "),
    @(42, "6.53",       "  +1.69%  "),
    @(43, "42.69",      "  +1.57%  "),
    @(44, "0.0316",     "  +1.55%  "),
    @(45, "25.98",      "  +11.83%  "),
    @(46, "0.772",      "  +1.86%  "),
    @(49, "308.47",     "  +5.30%  "),
    @(50, "0.870",      "  +4.48%  "),
    @(51, $null,        "  -1.30%  ")
)

foreach ($u in $updates) {
    $row = $u[0]
    $newD = $u[1]
    $newE = $u[2]
    if ($null -ne $newD) {
        Set-TextCell $ws.Cells.Item($row, 4) $newD
    }
    if ($null -ne $newE) {
        $ws.Cells.Item($row, 5).Value = $newE
    }
}

# Rows 47 and 48 swap their coin content (ONDO <-> Stellar) with new values.
$ws.Cells.Item(47, 2).Value = "Stellar"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell $ws.Cells.Item(47, 4) "0.112"
$ws.Cells.Item(47, 5).Value = "  +8.89%  "

$ws.Cells.Item(48, 2).Value = "ONDO"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-TextCell $ws.Cells.Item(48, 4) "1.08"
$ws.Cells.Item(48, 5).Value = "  +1.92%  "
